$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "UniformF"

# Add new row 16, mirroring the pattern of row 15 (A:M)
# Copy formatting from row 15 so the new row matches existing styling (bold/border on col A)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

for ($col = 3; $col -le 13; $col++) {
    $ws.Cells.Item(16, $col).Value = 1
}
